$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Joseline contributed additional screens: New Bike Part Order, New Product Manufacturer, New Bike Manufacturer ---
$ws.Range("B31").Value = "Joseline"
$ws.Range("B34").Value = "Joseline"

# Row 40: rename "New Bike Manufacturer" -> "newBikeManufacturer", contributor Joseline
$ws.Range("A40").Value = "newBikeManufacturer"
$ws.Range("B40").Value = "Joseline"

# --- New "Notes" content for Joseline's newly reviewed screens (rows 24 & 25) ---
$ws.Range("C24").Value = "need to adjust - remove table just text fields/labels"
$ws.Range("C25").Value = "need to adjust - remove table just text fields/labels"

# De-highlight rows 24 & 25 (no longer "main" screens) to match the other sub-screen rows
$ws.Range("A33").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A33").Copy()
$ws.Range("A25").PasteSpecial(-4122)

# Highlight rows 31, 34 & 40 like the other "main" screen rows
$ws.Range("A14").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("A40").PasteSpecial(-4122)

# Update the last-active-cell selection
$ws.Range("C13").Select()
